# Add the "LTSD Parameters" block (activity-detection parameters) to the
# TS40 worksheet: replace the stray "Test" label in K2 with the real
# section header, then populate the Right/Left threshold & window table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section header (was "Test")
$ws.Range("K2").Value = "LTSD Parameters"

# Right / Left sub-headers
$ws.Range("K3").Value = "Right"
$ws.Range("M3").Value = "Left"

# Column labels for each side
$ws.Range("K4").Value = "Threshols"
$ws.Range("L4").Value = "Win"
$ws.Range("M4").Value = "Threshold"
$ws.Range("N4").Value = "Win"

# Parameter values - stored as text (e.g. "6.0") rather than numbers, so
# force a text number format before assigning, then drop the formatting
# again so the cells keep the workbook's default (unstyled) look.
# Order matches the original authoring order (Right threshold/win, then
# Left threshold/win) so the shared-string table comes out the same way.
$paramCells = "M5", "K5", "L5", "N5"
$paramValues = "6.0", "5.5", "100.0", "300.0"
for ($i = 0; $i -lt $paramCells.Length; $i++) {
    $cell = $ws.Range($paramCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $paramValues[$i]
    $cell.ClearFormats()
}

# Move the active selection to N5, matching the saved view state
$ws.Range("N5").Select() | Out-Null
